$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update measured distance values (D column) ---
$ws.Range("D12").Value2 = 265
$ws.Range("D13").Value2 = 237
$ws.Range("D15").Value2 = 183
$ws.Range("D16").Value2 = 170
$ws.Range("D17").Value2 = 150

# --- Update the slope formula (I2:I16) to reference the next row instead of row+5 ---
$ws.Range("I2").Formula = "=(E3-E2)/(D3-D2)"
$ws.Range("I3:I16").Formula = "=(E4-E3)/(D4-D3)"

# --- Move/resize the embedded chart ---
$co = $ws.ChartObjects().Item(1)
$co.Left = 394.125
$co.Top = 72.75
$co.Width = 799.375
$co.Height = 334.5

# --- Update the selected cell/range shown when the sheet is active ---
$ws.Range("D18").Select() | Out-Null
